# Fruta / hortaliza, semanal
# A new weekly price-report row is inserted into the "Pomelo" data table,
# pushing the existing rows 36-135 down by one (to 37-136) and growing the
# used range from A1:T135 to A1:T136. The new row 36 carries the latest
# "Start Ruby" / "Primera" observation for Vega Modelo de Temuco.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 36; Excel shifts rows 36:135 down to 37:136 and
# extends the sheet dimension to A1:T136 automatically.
$ws.Rows.Item(36).Insert()

# Populate the newly inserted row 36 with the new weekly record.
$ws.Cells.Item(36, 1).Value  = 10
$ws.Cells.Item(36, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(36, 3).Value  = "La Araucanía"
$ws.Cells.Item(36, 4).Value  = 44459
$ws.Cells.Item(36, 5).Value  = 9
$ws.Cells.Item(36, 6).Value  = "Fruta"
$ws.Cells.Item(36, 7).Value  = 100102
$ws.Cells.Item(36, 8).Value  = "Cítricos"
$ws.Cells.Item(36, 9).Value  = 100102006
$ws.Cells.Item(36, 10).Value = "Pomelo"
$ws.Cells.Item(36, 11).Value = "Start Ruby"
$ws.Cells.Item(36, 12).Value = "Primera"
$ws.Cells.Item(36, 13).Value = 80
$ws.Cells.Item(36, 14).Value = 10000
$ws.Cells.Item(36, 15).Value = 10000
$ws.Cells.Item(36, 16).Value = 10000
$ws.Cells.Item(36, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(36, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(36, 19).Value = 667
$ws.Cells.Item(36, 20).Value = 15
